$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New cell data for columns H:M (HW resources info), rows 1-18 and 23-26.
# ---------------------------------------------------------------------------
$data = @(
    @(1, 8, 'SCH signal'),
    @(1, 9, 'SCH N'),
    @(1, 10, 'Jetson Socket N'),
    @(1, 11, 'DevKit Signal'),
    @(1, 12, 'SOC signal'),
    @(1, 13, 'Real HW'),
    @(2, 8, 'AUDIO_MCLK'),
    @(2, 9, 1),
    @(2, 10, 7),
    @(2, 11, 'AUDIO_I2S_MCLK_3V3'),
    @(2, 12, 'AUDIO_MCLK'),
    @(2, 13, 'AUD_MCLK'),
    @(3, 8, 'I2S_SCLK'),
    @(3, 9, 2),
    @(3, 10, 12),
    @(3, 11, 'AUDIO_I2S_SRCLK_3V3'),
    @(3, 12, 'I2S0_SCLK'),
    @(3, 13, 'I2S-1'),
    @(4, 8, 'O20'),
    @(4, 9, 3),
    @(4, 10, 13),
    @(4, 11, 'AUDIO_CDC_IRQ_LVL'),
    @(4, 12, 'GPIO20'),
    @(4, 13, 'GPIO3_PJ.05'),
    @(5, 8, 'O16'),
    @(5, 9, 4),
    @(5, 10, 18),
    @(5, 11, 'MDM_WAKE_AP_LVL'),
    @(5, 12, 'GPIO16'),
    @(5, 13, 'GPIO3_PY.01'),
    @(6, 8, 'SPI_MOSI'),
    @(6, 9, 5),
    @(6, 10, 19),
    @(6, 11, 'SPI1_MOSI_3V3'),
    @(6, 12, 'SPI1_MOSI'),
    @(6, 13, 'SPI-4'),
    @(7, 8, 'SPI_MISO'),
    @(7, 9, 6),
    @(7, 10, 21),
    @(7, 11, 'SPI1_MISO_3V3'),
    @(7, 12, 'SPI1_MISO'),
    @(8, 8, 'SPI_CLK'),
    @(8, 9, 7),
    @(8, 10, 23),
    @(8, 11, 'SPI1_SCK_3V3'),
    @(8, 12, 'SPI1_CLK'),
    @(9, 8, 'SPI_CS'),
    @(9, 9, 8),
    @(9, 10, 24),
    @(9, 11, 'SPI1_CS0_3V3'),
    @(9, 12, 'SPI1_CS0#'),
    @(10, 8, 'I2C_SDA'),
    @(10, 9, 9),
    @(10, 10, 27),
    @(10, 11, 'I2C_GP1_DAT_3V3'),
    @(10, 12, 'I2C_GP1_DAT'),
    @(10, 13, 'I2C-1'),
    @(11, 8, 'I2C_SCL'),
    @(11, 9, 10),
    @(11, 10, 28),
    @(11, 11, 'I2C_GP1_CLK_3V3'),
    @(11, 12, 'I2C_GP1_CLK'),
    @(12, 8, 'O19'),
    @(12, 9, 11),
    @(12, 10, 29),
    @(12, 11, 'AUD_RST_LVL'),
    @(12, 12, 'GPIO19'),
    @(12, 13, 'GPIO3_PJ.06'),
    @(13, 8, 'O9'),
    @(13, 9, 12),
    @(13, 10, 31),
    @(13, 11, 'MOTION_INT_AP_L_LVL'),
    @(13, 12, 'GPIO9'),
    @(13, 13, 'GPIO3_PAA.02'),
    @(14, 8, 'O11'),
    @(14, 9, 13),
    @(14, 10, 33),
    @(14, 11, 'AP_WAKE_BT_3V3'),
    @(14, 12, 'GPIO11'),
    @(14, 13, 'GPIO3_PI.05'),
    @(15, 8, 'I2S_LRCLK'),
    @(15, 9, 14),
    @(15, 10, 35),
    @(15, 11, 'AUDIO_I2S_SFSYNC_3V3'),
    @(15, 12, 'I2S0_LRCLK'),
    @(15, 13, 'I2S-1'),
    @(16, 8, 'O8'),
    @(16, 9, 15),
    @(16, 10, 37),
    @(16, 11, 'SAR_TOUT_LVL'),
    @(16, 12, 'GPIO8'),
    @(16, 13, 'GPIO3_PI.04'),
    @(17, 8, 'I2S_SDIN'),
    @(17, 9, 16),
    @(17, 10, 38),
    @(17, 11, 'AUDIO_I2S_SIN_3V3'),
    @(17, 12, 'I2S0_SDIN'),
    @(17, 13, 'I2S-1'),
    @(18, 8, 'I2S_SDOUT'),
    @(18, 9, 17),
    @(18, 10, 40),
    @(18, 11, 'AUDIO_I2S_SOUT_3V3'),
    @(18, 12, 'I2S0_SDOUT'),
    @(23, 10, 23),
    @(23, 11, 'I2S1_CLK'),
    @(23, 12, 'I2S1_CLK'),
    @(23, 13, 'I2S-2'),
    @(24, 10, 24),
    @(24, 11, 'I2S1_SDOUT'),
    @(24, 12, 'I2S1_SDOUT'),
    @(25, 10, 25),
    @(25, 11, 'I2S1_SDIN'),
    @(25, 12, 'I2S1_SDIN'),
    @(26, 10, 26),
    @(26, 11, 'I2S1_LRCLK'),
    @(26, 12, 'I2S1_LRCLK')
)

# ---------------------------------------------------------------------------
# 2. Formatting -- applied to ranges FIRST (in the order below) so that the
#    resulting cellXfs/fonts come out in a predictable sequence, then values
#    are written on top without disturbing the formats.
# ---------------------------------------------------------------------------

# Body font (size 14, regular) with no explicit alignment -> H2:L18
$bodyPlain = $ws.Range("H2:L18")
$bodyPlain.Font.Size = 14
$bodyPlain.Font.Bold = $false

# Header first cell (size 14, bold, left aligned) -> H1
$headerFirst = $ws.Range("H1")
$headerFirst.Font.Size = 14
$headerFirst.Font.Bold = $true
$headerFirst.HorizontalAlignment = -4131

# Rest of header row (size 14, bold, default alignment) -> I1:M1
$headerRest = $ws.Range("I1:M1")
$headerRest.Font.Size = 14
$headerRest.Font.Bold = $true

# "Real HW" column cells that are centered only -> M3, M15
$centerOnly = $ws.Range("M3")
$centerOnly.Font.Size = 14
$centerOnly.HorizontalAlignment = -4108
$centerOnly2 = $ws.Range("M15")
$centerOnly2.Font.Size = 14
$centerOnly2.HorizontalAlignment = -4108

# "Real HW" merged groups -> centered + vertically centered
$mergedGroups = @("M6:M9", "M10:M11", "M17:M18", "M23:M26")
foreach ($addr in $mergedGroups) {
    $rng = $ws.Range($addr)
    $rng.Font.Size = 14
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4108
}

# Remaining plain "Real HW" column cells (size 14, no alignment) -> same
# style as H2:L18 (reused automatically since it is identical formatting)
$plainM = @("M2", "M4", "M5", "M12", "M13", "M14", "M16")
foreach ($addr in $plainM) {
    $rng = $ws.Range($addr)
    $rng.Font.Size = 14
    $rng.Font.Bold = $false
}

# ---------------------------------------------------------------------------
# 3. Write the values now that every touched cell already has its format.
# ---------------------------------------------------------------------------
foreach ($item in $data) {
    $r = $item[0]
    $c = $item[1]
    $v = $item[2]
    $ws.Cells.Item($r, $c).Value = $v
}

# ---------------------------------------------------------------------------
# 4. Merge the "Real HW" cells that share one value across several rows.
# ---------------------------------------------------------------------------
$ws.Range("M23:M26").Merge() | Out-Null
$ws.Range("M17:M18").Merge() | Out-Null
$ws.Range("M6:M9").Merge() | Out-Null
$ws.Range("M10:M11").Merge() | Out-Null

# ---------------------------------------------------------------------------
# 5. Row heights.
# ---------------------------------------------------------------------------
for ($r = 1; $r -le 18; $r++) {
    $ws.Rows.Item($r).RowHeight = 18.5
}
$ws.Rows.Item(20).RowHeight = 15
$ws.Rows.Item(23).RowHeight = 18.5
$ws.Rows.Item(24).RowHeight = 18.5

# ---------------------------------------------------------------------------
# 6. Column widths for the new columns (best-fit approximations).
# ---------------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 14.5
$ws.Columns.Item(10).ColumnWidth = 16.67
$ws.Columns.Item(11).ColumnWidth = 25.83
$ws.Columns.Item(12).ColumnWidth = 18.83
$ws.Columns.Item(13).ColumnWidth = 15.5

# ---------------------------------------------------------------------------
# 7. Misc sheet-level tweaks.
# ---------------------------------------------------------------------------
$ws.Range("M16").Select() | Out-Null
